$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1219.3914
$ws.Range("I19").Value = 1213.8182
$ws.Range("K19").Value = 1213.8182
$ws.Range("M19").Value = -1038.8182
$ws.Range("H28").Value = 8023.385
$ws.Range("I28").Value = 8023.385
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 8023.385
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -7538.385
$ws.Range("N28").ClearContents()
$ws.Range("H76").Value = 2937.6
$ws.Range("I76").Value = 2972.5
$ws.Range("K76").Value = 2972.5
$ws.Range("M76").Value = -2657.5
$ws.Range("H79").Value = 2937.6
$ws.Range("I79").Value = 2972.5
$ws.Range("K79").Value = 2972.5
$ws.Range("M79").Value = -1880.5
$ws.Range("H124").Value = 49800
$ws.Range("J124").Value = 49800
$ws.Range("L124").Value = 49800
$ws.Range("N124").Value = -59620
$ws.Range("H137").Value = 1448.0416
$ws.Range("J137").Value = 2361.3333
$ws.Range("L137").Value = 7083.999899999999
$ws.Range("N137").Value = -12183.9999
$ws.Range("H138").Value = 1877.08
$ws.Range("I138").Value = 663.7059
$ws.Range("J138").Value = 2125.6023
$ws.Range("K138").Value = 1991.1177
$ws.Range("L138").Value = 6376.8069
$ws.Range("M138").Value = 3148.8823
$ws.Range("N138").Value = -16656.8069
$ws.Range("H141").Value = 11079.454
$ws.Range("I141").Value = 12685.444
$ws.Range("J141").Value = 3852.5
$ws.Range("K141").Value = 38056.33199999999
$ws.Range("L141").Value = 11557.5
$ws.Range("M141").Value = -32876.33199999999
$ws.Range("N141").Value = -21917.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3127.75
$ws.Range("I32").Value = 3147.2708
$ws.Range("K32").Value = 3147.2708
$ws.Range("M32").Value = -2860.2708
$ws.Range("H61").Value = 1396.5714
$ws.Range("I61").Value = 1021.7778
$ws.Range("K61").Value = 1021.7778
$ws.Range("M61").Value = -809.7778
$ws.Range("H74").Value = 1203.421
$ws.Range("I74").Value = 1079.7693
$ws.Range("J74").Value = 1471.3334
$ws.Range("K74").Value = 1079.7693
$ws.Range("L74").Value = 1471.3334
$ws.Range("M74").Value = -205.7692999999999
$ws.Range("N74").Value = -3219.3334
$ws.Range("H77").Value = 1203.421
$ws.Range("I77").Value = 1079.7693
$ws.Range("J77").Value = 1471.3334
$ws.Range("K77").Value = 5398.8465
$ws.Range("L77").Value = 7356.666999999999
$ws.Range("M77").Value = -1030.8465
$ws.Range("N77").Value = -16092.667
$ws.Range("H122").Value = 1098.7333
$ws.Range("I122").Value = 1098.7333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3296.199900000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -846.1999000000005
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2220.6511
$ws.Range("I132").Value = 1959
$ws.Range("J132").Value = 3365.375
$ws.Range("K132").Value = 5877
$ws.Range("L132").Value = 10096.125
$ws.Range("M132").Value = -3347
$ws.Range("N132").Value = -15156.125
$ws.Range("H136").Value = 1396.5714
$ws.Range("I136").Value = 1021.7778
$ws.Range("K136").Value = 3065.3334
$ws.Range("M136").Value = -515.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 142859630
$ws.Range("I105").Value = 166669220
$ws.Range("K105").Value = 166669220
$ws.Range("M105").Value = -166667473
$ws.Range("H134").Value = 3977.8948
$ws.Range("I134").Value = 1355.0312
$ws.Range("K134").Value = 4065.0936
$ws.Range("M134").Value = -1530.0936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 10050
$ws.Range("I23").Value = 4200
$ws.Range("K23").Value = 4200
$ws.Range("M23").Value = -3960
$ws.Range("H27").Value = 10050
$ws.Range("I27").Value = 4200
$ws.Range("K27").Value = 4200
$ws.Range("M27").Value = -4008
$ws.Range("H31").Value = 1487.3448
$ws.Range("I31").Value = 938.5333
$ws.Range("K31").Value = 938.5333
$ws.Range("M31").Value = -643.5333
$ws.Range("H34").Value = 1487.3448
$ws.Range("I34").Value = 938.5333
$ws.Range("K34").Value = 938.5333
$ws.Range("M34").Value = -736.5333
$ws.Range("H99").Value = 1705.5454
$ws.Range("I99").Value = 1858
$ws.Range("K99").Value = 1858
$ws.Range("M99").Value = -360
$ws.Range("H126").Value = 1705.5454
$ws.Range("I126").Value = 1858
$ws.Range("K126").Value = 5574
$ws.Range("M126").Value = -3104
$ws.Range("H132").Value = 5224.294
$ws.Range("I132").Value = 5932.64
$ws.Range("J132").Value = 3256.6667
$ws.Range("K132").Value = 17797.92
$ws.Range("L132").Value = 9770.000100000001
$ws.Range("M132").Value = -15267.92
$ws.Range("N132").Value = -14830.0001
$ws.Range("H134").Value = 8334544.5
$ws.Range("I134").Value = 10753855
$ws.Range("J134").Value = 1365.2222
$ws.Range("K134").Value = 32261565
$ws.Range("L134").Value = 4095.6666
$ws.Range("M134").Value = -32259030
$ws.Range("N134").Value = -9165.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 387
$ws.Range("I20").Value = 387
$ws.Range("K20").Value = 1161
$ws.Range("M20").Value = -934
$ws.Range("H22").Value = 2744.4443
$ws.Range("J22").Value = 2785.7144
$ws.Range("L22").Value = 8357.143199999999
$ws.Range("N22").Value = -8695.143199999999
$ws.Range("H27").Value = 2744.4443
$ws.Range("J27").Value = 2785.7144
$ws.Range("L27").Value = 8357.143199999999
$ws.Range("N27").Value = -8561.143199999999
$ws.Range("H68").Value = 1744.8214
$ws.Range("I68").Value = 747.06665
$ws.Range("J68").Value = 2896.077
$ws.Range("K68").Value = 2241.19995
$ws.Range("L68").Value = 8688.231
$ws.Range("M68").Value = -1430.19995
$ws.Range("N68").Value = -10310.231
$ws.Range("H71").Value = 1744.8214
$ws.Range("I71").Value = 747.06665
$ws.Range("J71").Value = 2896.077
$ws.Range("K71").Value = 6723.59985
$ws.Range("L71").Value = 26064.693
$ws.Range("M71").Value = -2667.59985
$ws.Range("N71").Value = -34176.693
$ws.Range("H107").Value = 3458.543
$ws.Range("I107").Value = 532.375
$ws.Range("J107").Value = 9842.909
$ws.Range("K107").Value = 1597.125
$ws.Range("L107").Value = 29528.727
$ws.Range("M107").Value = 322.875
$ws.Range("N107").Value = -33368.727
$ws.Range("H131").Value = 22728758
$ws.Range("J131").Value = 1586.8158
$ws.Range("L131").Value = 4760.4474
$ws.Range("N131").Value = -14840.4474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1425.1666
$ws.Range("I68").Value = 1113
$ws.Range("K68").Value = 1113
$ws.Range("M68").Value = -364
$ws.Range("H71").Value = 1425.1666
$ws.Range("I71").Value = 1113
$ws.Range("K71").Value = 5565
$ws.Range("M71").Value = -1821
$ws.Range("H100").Value = 2200.8
$ws.Range("I100").Value = 1750
$ws.Range("K100").Value = 1750
$ws.Range("M100").Value = -1209
$ws.Range("H132").Value = 45791.74
$ws.Range("I132").Value = 2015.6154
$ws.Range("J132").Value = 102700.7
$ws.Range("K132").Value = 6046.8462
$ws.Range("L132").Value = 308102.1
$ws.Range("M132").Value = -3516.8462
$ws.Range("N132").Value = -313162.1
$ws.Range("H136").Value = 7616.375
$ws.Range("I136").Value = 12138.777
$ws.Range("J136").Value = 1801.8572
$ws.Range("K136").Value = 36416.331
$ws.Range("L136").Value = 5405.571599999999
$ws.Range("M136").Value = -33866.331
$ws.Range("N136").Value = -10505.5716
